# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "last row" date format (currently on A90) before it moves.
$lastRowDateFormat = $ws.Range("A90").NumberFormat
$normalDateFormat = $ws.Range("A2").NumberFormat

# The old last row (row 90) becomes a regular row, so it takes on the
# normal date format used by all the other non-last rows.
$ws.Range("A90").NumberFormat = $normalDateFormat

# Append the new day's data as the new last row (row 91).
$ws.Range("A91").Value = 45678
$ws.Range("B91").Value = 213
$ws.Range("C91").Value = 211
$ws.Range("D91").Value = 213

# The newly appended row becomes the new last row, so it takes on the
# "last row" date format that row 90 used to have.
$ws.Range("A91").NumberFormat = $lastRowDateFormat
